# All ADC channels have been tested and are functional.
# Fill in the newly-verified engine_config sensor rows (17-20): byte counts
# in column A and descriptions in column D. The existing shared formulas in
# columns B/C (address bookkeeping) recompute automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A17").Value = 1
$ws.Range("D17").Value = "Throttle position sensor Low (engine_config) 8 bits "

$ws.Range("A18").Value = 1
$ws.Range("D18").Value = "Throttle position sensor high (engine_config) 8 bits "

$ws.Range("A19").Value = 1

$ws.Range("A20").Value = 2
# Write D20 before D19 so the new shared-string table entries land in the
# same order as the source workbook (Manifold high before Manifold Low).
$ws.Range("D20").Value = "Manifold absolute pressure high (engine_config) 8 bits "
$ws.Range("D19").Value = "Manifold absolute pressure Low (engine_config) 8 bits "

# Move/restore the active selection to A18, matching the saved view state.
$ws.Range("A18").Select() | Out-Null
